$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Components that have now actually been purchased: record the price paid
# for the ESP32 module and mark it (and the two resistors) as "Bought?" = y
$ws.Range("D4").Value = 7
$ws.Range("F4").Value = "y"
$ws.Range("F14").Value = "y"
$ws.Range("F15").Value = "y"

# Leave the selection where the author left it when they saved the file
$ws.Range("G24").Select() | Out-Null
